$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("B2").Value = "Quantiferon TB Gold (Modified)"
$ws.Range("C2").Value = 4001158
$ws.Range("D2").Value = 11833232
$ws.Range("E2").Value = "Negative"

# Update row 3
$ws.Range("B3").Value = "Quantiferon TB Gold (Modified)"
$ws.Range("C3").Value = 4001033
$ws.Range("D3").Value = 11826881
$ws.Range("E3").Value = "Positive"

# Delete rows 4 through 9 (no longer present in the data)
$ws.Range("A4:F9").Delete()
